$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3: "BurnDown Chart" title -> "Burndown Chart" (typo fix + run split)
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$titleShape3 = $s3.Shapes.Item(1)
$tr3 = $titleShape3.TextFrame.TextRange

$f = $tr3.Find("BurnDown", 0)
$f.Text = "Burndown"

$anchor = $tr3.Find(" Chart" + [char]9, 0)
$anchor.InsertAfter("Chart" + [char]9)
$anchor.Text = " "

# ---------------------------------------------------------------------------
# Slide 9: "Oportunidades de mejora (1)" - fix wording
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$bodyShape9 = $s9.Shapes.Item(2)
$tr9 = $bodyShape9.TextFrame.TextRange

$f = $tr9.Find("Objetivo", 0)
$f.Text = "Reuniones"

$f = $tr9.Find(" reunions ", 0)
$f.Text = " - "

$f = $tr9.Find("planeamiento", 0)
$f.Text = "Planeamiento"

# Append " aprobado" right after the existing "definido" run
$f = $tr9.Find("definido", 0)
$f.InsertAfter("aprobado")
$f.InsertAfter(" ")
$f.Text = "definido"

# ---------------------------------------------------------------------------
# Slide 10: "Oportunidades de mejora (2)" - fix wording
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$bodyShape10 = $s10.Shapes.Item(2)
$tr10 = $bodyShape10.TextFrame.TextRange

$f = $tr10.Find("separados y se ", 0)
$f.Text = "separados y "

$f = $tr10.Find("bloqueados", 0)
$f.Text = "bloqueados."

$f = $tr10.Find("reporte", 0)
$f.Text = "informe"

# Rebuild ", el cliente no llegaba a consumir todos." into the new sentence
$anchor = $tr10.Find(", el ", 0)
$anchor.InsertAfter(" al ")
$anchor.InsertAfter("agobiar")
$anchor.InsertAfter("para no ")
$anchor.InsertAfter(" ")
$anchor.InsertAfter("internos")
$anchor.InsertAfter(" ")
$anchor.InsertAfter("detalles")
$anchor.InsertAfter(" a ")
$anchor.InsertAfter("tanto")
$anchor.InsertAfter(" y no ")
$anchor.InsertAfter("completas")
$anchor.InsertAfter(" ")
$anchor.InsertAfter("funcionalidades")
$anchor.InsertAfter(" a ")
$anchor.InsertAfter("orientaci" + [char]0x00F3 + "n")
$anchor.InsertAfter("con ")
$anchor.Text = ", "

$f = $tr10.Find(" no llegaba a consumir todos", 0)
$f.Text = ""
